# Enabled_Controller_BOM: split single BOM sheet into "USB" and "Wireless" variants.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- 1. Fix up the original sheet's data first (USB variant), while it is
#        still the only sheet, then clone it for the Wireless variant so
#        both sheets share headers / unrelated rows without retyping them.

# Row 2: swap the old "ITSY BITSY 32U4" line for the new USB microcontroller.
$ws1.Range("H2").Value = "https://www.digikey.ca/en/products/detail/adafruit-industries-llc/3727/8346575"
$ws1.Range("C2").Value = "1528-2554-ND"
$ws1.Range("E2").Value = "ITSYBITSY M0 EXPRESS ATSAMD21"
$ws1.Range("D2").Value = 3727
$ws1.Range("F2").Value = 16.91
$ws1.Range("G2").Value = 16.91
$ws1.Rows.Item(2).AutoFit()

$h2 = $null
foreach ($hl in $ws1.Hyperlinks) {
    if ($hl.Range.Address() -eq '$H$2') { $h2 = $hl }
}
if ($h2 -ne $null) {
    $h2.Address = "https://www.digikey.ca/en/products/detail/adafruit-industries-llc/3727/8346575"
}

# Row 7 / Row 8: fill in the previously-blank "Total Price" column.
$ws1.Range("G7").Value = 0.15
$ws1.Range("G8").Value = 0.15

# Row 11: total no longer includes the (now removed) row 9 blank line.
$ws1.Range("G11").Formula = "=F2*B2+F3*B3+F4*B4+F5*B5+F6*B6+F7*B7+F8*B8"

# --- 2. Clone the now-updated USB sheet to become the Wireless sheet,
#        placed immediately after it.
$ws1.Copy($null, $ws1)
$ws2 = $wb.Worksheets.Item(2)

# --- 3. Rename both sheets.
$ws1.Name = "USB Enabled_Controller_BOM"
$ws2.Name = "Wireless Enabled_Controller_BOM"

# --- 4. Wireless sheet gets its own microcontroller in row 2.
$ws2.Range("C2").Value = "1528-4481-ND"
$ws2.Range("E2").Value = "`tITSYBITSY NRF52840 EXPRESS BLE"
$ws2.Range("H2").Value = "https://www.digikey.ca/en/products/detail/adafruit-industries-llc/4481/11497502"
$ws2.Range("D2").Value = 4481
$ws2.Range("F2").Value = 25.4
$ws2.Range("G2").Value = 25.4
$ws2.Rows.Item(2).AutoFit()

$h2b = $null
foreach ($hl in $ws2.Hyperlinks) {
    if ($hl.Range.Address() -eq '$H$2') { $h2b = $hl }
}
if ($h2b -ne $null) {
    $h2b.Address = "https://www.digikey.ca/en/products/detail/adafruit-industries-llc/4481/11497502"
}

# --- 5. Trim the Wireless sheet down to the 8-column / 11-row BOM (no
#        trailing blank I/J columns or rows 10, 12-14 like the source sheet
#        used to have).
$ws2.Range("I1:J14").Clear()
$ws2.Range("A10:H10").Clear()
$ws2.Range("A12:H14").Clear()
$ws2.Range("A11:E11").Clear()
$ws2.Range("H11").Clear()

# --- 6. Selections: USB sheet no longer the active tab; Wireless does.
$ws1.Range("A1:H11").Select()
$ws2.Range("H2").Select()
$ws2.Activate()
